# Added multiple sheet handling capability
# - keep existing "cars" sheet (sheet1) data untouched
# - add two more sheets ("cars2", "cars3") with car data, positioned after "cars"
# - update view/selection state to match the edited workbook

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- add the two new sheets, right after "cars" and in order ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "cars2"

$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "cars3"

# --- cars2: same shape/values as "cars", but header "zcarModel" ---
$ws2.Cells.Item(1, 1).Value = "zcarModel"
$ws2.Cells.Item(1, 2).Value = "price"
$ws2.Cells.Item(1, 3).Value = "color"

$ws2.Cells.Item(2, 1).Value = "Audi"
$ws2.Cells.Item(2, 2).Value = 10000
$ws2.Cells.Item(2, 3).Value = "blue"

$ws2.Cells.Item(3, 1).Value = "BMW"
$ws2.Cells.Item(3, 2).Value = 15000
$ws2.Cells.Item(3, 3).Value = "red"

$ws2.Cells.Item(4, 1).Value = "Mercedes"
$ws2.Cells.Item(4, 2).Value = 20000
$ws2.Cells.Item(4, 3).Value = "yellow"

$ws2.Cells.Item(5, 1).Value = "Porsche"
$ws2.Cells.Item(5, 2).Value = 30000
$ws2.Cells.Item(5, 3).Value = "green"

# --- cars3: new data set, header "ycarModel" ---
$ws3.Cells.Item(1, 1).Value = "ycarModel"
$ws3.Cells.Item(1, 2).Value = "price"
$ws3.Cells.Item(1, 3).Value = "color"

$ws3.Cells.Item(2, 1).Value = "Chevy"
$ws3.Cells.Item(2, 2).Value = 15000
$ws3.Cells.Item(2, 3).Value = "white"

$ws3.Cells.Item(3, 1).Value = "Mazda"
$ws3.Cells.Item(3, 2).Value = 19000
$ws3.Cells.Item(3, 3).Value = "silver"

$ws3.Cells.Item(4, 1).Value = "Ford"
$ws3.Cells.Item(4, 2).Value = 20000
$ws3.Cells.Item(4, 3).Value = "red"

$ws3.Cells.Item(5, 1).Value = "Toyota"
$ws3.Cells.Item(5, 2).Value = 16000
$ws3.Cells.Item(5, 3).Value = "black"

# --- view/selection state ---
# sheet1 ("cars"): selection covers the whole table, anchored visually at C5
[void]$ws1.Range("A1:C5").Select()

# sheet2 ("cars2"): selection parked on a single empty cell outside the data
[void]$ws2.Range("F11").Select()

# sheet3 ("cars3"): active sheet, selection parked on a single cell outside the data
[void]$ws3.Range("I3").Select()
[void]$ws3.Activate()

Write-Output "ok"
